# Update MN and GA
# On the MN sheet, the "gtu"/"ur"/"urine"/"urine_card" specimen-type rows
# (E7:E10) now map to a new "urinary" sti report pdf form field value
# instead of the previous "urine" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MN")

$ws.Range("E7:E10").Value = "urinary"

$ws.Activate()
$ws.Range("E12").Select()
